$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28: CCEP - Coca-Cola Europacific Partners
$ws.Range("C28").Value = "Consumer Staples"
$ws.Range("D28").Value = "Soft Drinks & Non-alcoholic Beverages"

# Row 41: DASH - DoorDash (fix company name capitalization + add sub-sector)
$ws.Range("B41").Value = "DoorDash"
$ws.Range("D41").Value = "Specialized Consumer Services"

# Row 69: MDB - MongoDB Inc.
$ws.Range("D69").Value = "Systems Software"
